# Highlight quantitative impact metrics (percentages, dollar amounts, large
# numbers) in bold + a dark slate color (#2C3E50) across the ATS long
# cartographic-professional resume for Dheeraj Chand.
#
# Word's Font.Color property takes a BGR-packed long (like the VBA RGB()
# macro), so #2C3E50 (R=0x2C,G=0x3E,B=0x50) must be passed as 0x503E2C.
$metricColor = 0x503E2C

$d = $word.ActiveDocument

# Bolds + colors the next (left-to-right, case-sensitive) occurrence of
# $metricText inside the paragraph $paraIndex (1-based, matching Word's
# Paragraphs collection), searching only after character offset $afterPos
# within that paragraph. Returns the end offset of the match so callers can
# chain multiple metrics in a single paragraph in left-to-right order.
function Highlight-Metric($paraIndex, $afterPos, $metricText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $paraEnd = $p.Range.End
    $startPos = $afterPos
    if ($startPos -eq $null -or $startPos -lt $p.Range.Start) {
        $startPos = $p.Range.Start
    }
    $rng = $d.Range($startPos, $paraEnd)
    $found = $rng.Find.Execute($metricText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Color = $metricColor
        return $rng.End
    }
    return $startPos
}

# Applies a sequence of metric substrings, in left-to-right order, within a
# single paragraph.
function Highlight-Metrics($paraIndex, $metrics) {
    $pos = $null
    foreach ($m in $metrics) {
        $pos = Highlight-Metric $paraIndex $pos $m
    }
}

# Partner - Siege Analytics bullet: "... from 23% to 64%"
Highlight-Metrics 10 @("23%", "64%")

# Partner - Siege Analytics bullet: survey margin / turnout accuracy
Highlight-Metrics 12 @("±4.2%", "±2.1%", "71%", "87%")

# Partner - Siege Analytics bullet: boundary estimation cost savings
Highlight-Metrics 13 @("73.5%", "$4.7M")

# Partner - Siege Analytics bullet: FEC analysis valued over $2 trillion
Highlight-Metrics 14 @("$2")

# Data Products Manager bullet: ETL processing time reduction
Highlight-Metrics 24 @("57%")

# Key Achievements: revenue generation
Highlight-Metrics 50 @("$4.9M")

# Key Achievements: conversion rate improvement
Highlight-Metrics 51 @("23%")

# Key Achievements: platform impact analyst count
Highlight-Metrics 53 @("12,847")

Write-Output "Highlighted quantitative metrics across 8 bullet paragraphs."
